$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append five more days of case counts (missed days caught up on 9/27).
$ws.Range("A195").Value = "9/23/20"
$ws.Range("B195").Value = 28745

$ws.Range("A196").Value = "9/24/20"
$ws.Range("B196").Value = 28904

$ws.Range("A197").Value = "9/25/20"
$ws.Range("B197").Value = 29073

$ws.Range("A198").Value = "9/26/20"
$ws.Range("B198").Value = 29130

$ws.Range("A199").Value = "9/27/20"
$ws.Range("B199").Value = 29252

# Scroll the view down and move the selection to where editing left off.
$win = $excel.ActiveWindow
$win.ScrollRow = 173
$win.ScrollColumn = 1
[void]$ws.Range("C198").Select()
